# Refresh the crypto snapshot values (price + 1h volume change) and, for the
# three coins whose ranking swapped places this run, their name/link too.
# Values that read as plain numbers (no thousands-dot separator) are forced
# through as text so Excel doesn't strip significant leading/trailing zeros
# or flip tiny prices into scientific notation -- matching how the source
# data is stored (inline strings, not numeric cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '62.396.43'; AsText = $false },
    @{ Cell = 'E2'; Value = '  -2.59%  '; AsText = $false },
    @{ Cell = 'D3'; Value = '3.001.42'; AsText = $false },
    @{ Cell = 'E3'; Value = '  -4.55%  '; AsText = $false },
    @{ Cell = 'E4'; Value = '  +0.02%  '; AsText = $false },
    @{ Cell = 'D5'; Value = '553.96'; AsText = $true },
    @{ Cell = 'E5'; Value = '  -2.28%  '; AsText = $false },
    @{ Cell = 'D6'; Value = '151.96'; AsText = $true },
    @{ Cell = 'E6'; Value = '  -7.08%  '; AsText = $false },
    @{ Cell = 'E7'; Value = '  +0.00%  '; AsText = $false },
    @{ Cell = 'D8'; Value = '0.570'; AsText = $true },
    @{ Cell = 'E8'; Value = '  -2.52%  '; AsText = $false },
    @{ Cell = 'D9'; Value = '3.004.92'; AsText = $false },
    @{ Cell = 'E9'; Value = '  -4.19%  '; AsText = $false },
    @{ Cell = 'D10'; Value = '0.113'; AsText = $true },
    @{ Cell = 'E10'; Value = '  -2.60%  '; AsText = $false },
    @{ Cell = 'D11'; Value = '6.31'; AsText = $true },
    @{ Cell = 'E11'; Value = '  -5.71%  '; AsText = $false },
    @{ Cell = 'D12'; Value = '0.367'; AsText = $true },
    @{ Cell = 'E12'; Value = '  -3.55%  '; AsText = $false },
    @{ Cell = 'D13'; Value = '3.529.14'; AsText = $false },
    @{ Cell = 'E13'; Value = '  -4.21%  '; AsText = $false },
    @{ Cell = 'E14'; Value = '  -3.17%  '; AsText = $false },
    @{ Cell = 'D15'; Value = '62.547.30'; AsText = $false },
    @{ Cell = 'E15'; Value = '  -2.03%  '; AsText = $false },
    @{ Cell = 'D16'; Value = '23.75'; AsText = $true },
    @{ Cell = 'E16'; Value = '  -4.80%  '; AsText = $false },
    @{ Cell = 'B17'; Value = 'ShibaInu'; AsText = $false },
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; AsText = $false },
    @{ Cell = 'D17'; Value = '0.0000149'; AsText = $true },
    @{ Cell = 'E17'; Value = '  -3.18%  '; AsText = $false },
    @{ Cell = 'B18'; Value = 'WrappedEther'; AsText = $false },
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; AsText = $false },
    @{ Cell = 'D18'; Value = '3.004.15'; AsText = $false },
    @{ Cell = 'E18'; Value = '  -4.12%  '; AsText = $false },
    @{ Cell = 'D19'; Value = '394.64'; AsText = $true },
    @{ Cell = 'E19'; Value = '  -3.47%  '; AsText = $false },
    @{ Cell = 'D20'; Value = '5.12'; AsText = $true },
    @{ Cell = 'E20'; Value = '  -2.14%  '; AsText = $false },
    @{ Cell = 'D21'; Value = '11.98'; AsText = $true },
    @{ Cell = 'E21'; Value = '  -3.85%  '; AsText = $false },
    @{ Cell = 'D22'; Value = '6.69'; AsText = $true },
    @{ Cell = 'E22'; Value = '  -5.43%  '; AsText = $false },
    @{ Cell = 'D23'; Value = '1.00'; AsText = $true },
    @{ Cell = 'E23'; Value = '  -0.15%  '; AsText = $false },
    @{ Cell = 'D24'; Value = '65.19'; AsText = $true },
    @{ Cell = 'E24'; Value = '  -3.18%  '; AsText = $false },
    @{ Cell = 'D25'; Value = '0.468'; AsText = $true },
    @{ Cell = 'E25'; Value = '  -2.76%  '; AsText = $false },
    @{ Cell = 'E26'; Value = '  -6.68%  '; AsText = $false },
    @{ Cell = 'D27'; Value = '0.0₃0971'; AsText = $false },
    @{ Cell = 'E27'; Value = '  -4.97%  '; AsText = $false },
    @{ Cell = 'E28'; Value = '  -3.31%  '; AsText = $false },
    @{ Cell = 'D29'; Value = '0.998'; AsText = $true },
    @{ Cell = 'E29'; Value = '  -0.34%  '; AsText = $false },
    @{ Cell = 'E30'; Value = '  +0.00%  '; AsText = $false },
    @{ Cell = 'E31'; Value = '  -3.24%  '; AsText = $false },
    @{ Cell = 'D32'; Value = '20.58'; AsText = $true },
    @{ Cell = 'E32'; Value = '  -2.37%  '; AsText = $false },
    @{ Cell = 'D33'; Value = '160.12'; AsText = $true },
    @{ Cell = 'E33'; Value = '  +4.88%  '; AsText = $false },
    @{ Cell = 'D34'; Value = '4.69'; AsText = $true },
    @{ Cell = 'E34'; Value = '  -1.62%  '; AsText = $false },
    @{ Cell = 'D35'; Value = '6.04'; AsText = $true },
    @{ Cell = 'E35'; Value = '  -3.52%  '; AsText = $false },
    @{ Cell = 'D36'; Value = '1.09'; AsText = $true },
    @{ Cell = 'E36'; Value = '  -2.97%  '; AsText = $false },
    @{ Cell = 'E37'; Value = '  -2.56%  '; AsText = $false },
    @{ Cell = 'E38'; Value = '  -4.87%  '; AsText = $false },
    @{ Cell = 'D39'; Value = '2.459.44'; AsText = $false },
    @{ Cell = 'E39'; Value = '  -10.55%  '; AsText = $false },
    @{ Cell = 'D40'; Value = '3.93'; AsText = $true },
    @{ Cell = 'E40'; Value = '  -3.30%  '; AsText = $false },
    @{ Cell = 'B41'; Value = 'EnergySwap'; AsText = $false },
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; AsText = $false },
    @{ Cell = 'D41'; Value = '22.58'; AsText = $true },
    @{ Cell = 'E41'; Value = '  -3.68%  '; AsText = $false },
    @{ Cell = 'B42'; Value = 'OKB'; AsText = $false },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; AsText = $false },
    @{ Cell = 'D42'; Value = '37.52'; AsText = $true },
    @{ Cell = 'E42'; Value = '  -3.22%  '; AsText = $false },
    @{ Cell = 'D43'; Value = '0.662'; AsText = $true },
    @{ Cell = 'E43'; Value = '  -5.11%  '; AsText = $false },
    @{ Cell = 'D44'; Value = '0.0597'; AsText = $true },
    @{ Cell = 'E44'; Value = '  -4.44%  '; AsText = $false },
    @{ Cell = 'B45'; Value = 'VeChain'; AsText = $false },
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; AsText = $false },
    @{ Cell = 'D45'; Value = '0.0248'; AsText = $true },
    @{ Cell = 'E45'; Value = '  -3.71%  '; AsText = $false },
    @{ Cell = 'B46'; Value = 'FirstDigitalUSD'; AsText = $false },
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; AsText = $false },
    @{ Cell = 'D46'; Value = '0.999'; AsText = $true },
    @{ Cell = 'E46'; Value = '  +0.04%  '; AsText = $false },
    @{ Cell = 'D47'; Value = '4.95'; AsText = $true },
    @{ Cell = 'E47'; Value = '  -8.40%  '; AsText = $false },
    @{ Cell = 'D48'; Value = '0.0954'; AsText = $true },
    @{ Cell = 'E48'; Value = '  -2.17%  '; AsText = $false },
    @{ Cell = 'D49'; Value = '19.82'; AsText = $true },
    @{ Cell = 'E49'; Value = '  -4.86%  '; AsText = $false },
    @{ Cell = 'D50'; Value = '10.48'; AsText = $true },
    @{ Cell = 'E50'; Value = '  +0.36%  '; AsText = $false },
    @{ Cell = 'D51'; Value = '263.38'; AsText = $true },
    @{ Cell = 'E51'; Value = '  -7.20%  '; AsText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.AsText) {
        $origStyle = $cell.Style
        $cell.NumberFormat = '@'
        $cell.Value = $u.Value
        $cell.Style = $origStyle
    } else {
        $cell.Value = $u.Value
    }
}
